# "neue samples in die listen eingefuegt"
# Append 21 new sample rows (A/B value pairs) to Tabelle1, rows 309-329,
# each styled with the same thin black border used throughout the table,
# then leave the selection where Excel would land after typing the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-SampleRow {
    param($Row, $AValue, $BValue)

    $cellA = $ws.Cells.Item($Row, 1)
    $cellA.Value = $AValue
    $cellA.Borders.ColorIndex = 1
    $cellA.Borders.Weight = 2
    $cellA.Borders.LineStyle = 1

    $cellB = $ws.Cells.Item($Row, 2)
    $cellB.Value = $BValue
    $cellB.Borders.ColorIndex = 1
    $cellB.Borders.Weight = 2
    $cellB.Borders.LineStyle = 1
}

$newRows = @(
    @(309, 86, 55),
    @(310, 86, 58),
    @(311, 86, 65),
    @(312, 87, 1),
    @(313, 87, 31),
    @(314, 87, 33),
    @(315, 87, 61),
    @(316, 87, 62),
    @(317, 87, 63),
    @(318, 87, 64),
    @(319, 87, 65),
    @(320, 87, 70),
    @(321, 88, 1),
    @(322, 88, 31),
    @(323, 88, 33),
    @(324, 88, 61),
    @(325, 88, 62),
    @(326, 88, 63),
    @(327, 88, 64),
    @(328, 88, 65),
    @(329, 88, 70)
)

foreach ($row in $newRows) {
    Add-SampleRow $row[0] $row[1] $row[2]
}

# Matches the author's final selection/view after adding the samples.
$ws.Range("B330").Select()
